$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "26-JAN-26"
$ws.Range("B2").Value = "SM-438"
$ws.Range("C2").Value = "Nesma Airlines NE-153"
$ws.Range("D2").Value = 350
$ws.Range("E2").Value = 672
$ws.Range("F2").Value = -322
$ws.Range("G2").Value = 30
$ws.Range("H2").Value = 30
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = "LOW THREAT"
$ws.Range("K2").Value = "SAR"

# Row 3
$ws.Range("A3").Value = "26-JAN-26"
$ws.Range("B3").Value = "SM-438"
$ws.Range("C3").Value = "Nile Air NP-232"
$ws.Range("D3").Value = 354
$ws.Range("E3").Value = 672
$ws.Range("F3").Value = -318
$ws.Range("G3").Value = 30
$ws.Range("H3").Value = 30
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = "LOW THREAT"
$ws.Range("K3").Value = "SAR"

# Row 4
$ws.Range("A4").Value = "26-JAN-26"
$ws.Range("B4").Value = "SM-438"
$ws.Range("C4").Value = "flyadeal F3-911"
$ws.Range("D4").Value = 447
$ws.Range("E4").Value = 672
$ws.Range("F4").Value = -225
$ws.Range("G4").Value = 30
$ws.Range("H4").Value = 30
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = "LOW THREAT"
$ws.Range("K4").Value = "SAR"

# Row 5
$ws.Range("A5").Value = "26-JAN-26"
$ws.Range("B5").Value = "SM-438"
$ws.Range("C5").Value = "flynas XY-854"
$ws.Range("D5").Value = 529
$ws.Range("E5").Value = 672
$ws.Range("F5").Value = -143
$ws.Range("G5").Value = 40
$ws.Range("H5").Value = 30
$ws.Range("I5").Value = -10
$ws.Range("J5").Value = "LOW THREAT"
$ws.Range("K5").Value = "SAR"

# Row 6
$ws.Range("A6").Value = "26-JAN-26"
$ws.Range("B6").Value = "SM-438"
$ws.Range("C6").Value = "flynas XY-894"
$ws.Range("D6").Value = 599
$ws.Range("E6").Value = 672
$ws.Range("F6").Value = -73
$ws.Range("G6").Value = 40
$ws.Range("H6").Value = 30
$ws.Range("I6").Value = -10
$ws.Range("J6").Value = "LOW THREAT"
$ws.Range("K6").Value = "SAR"

# Row 7
$ws.Range("A7").Value = "26-JAN-26"
$ws.Range("B7").Value = "SM-438"
$ws.Range("C7").Value = "EgyptAir MS-688"
$ws.Range("D7").Value = 656
$ws.Range("E7").Value = 672
$ws.Range("F7").Value = -16
$ws.Range("G7").Value = 46
$ws.Range("H7").Value = 30
$ws.Range("I7").Value = -16
$ws.Range("J7").Value = "LOW THREAT"
$ws.Range("K7").Value = "SAR"

# Row 8
$ws.Range("A8").Value = "02-FEB-26"
$ws.Range("B8").Value = "SM-438"
$ws.Range("C8").Value = "flynas XY-854"
$ws.Range("D8").Value = 599
$ws.Range("E8").Value = 745
$ws.Range("F8").Value = -146
$ws.Range("G8").Value = 40
$ws.Range("H8").Value = 30
$ws.Range("I8").Value = -10
$ws.Range("J8").Value = "LOW THREAT"
$ws.Range("K8").Value = "SAR"

# Row 9
$ws.Range("A9").Value = "02-FEB-26"
$ws.Range("B9").Value = "SM-438"
$ws.Range("C9").Value = "flynas XY-894"
$ws.Range("D9").Value = 599
$ws.Range("E9").Value = 745
$ws.Range("F9").Value = -146
$ws.Range("G9").Value = 40
$ws.Range("H9").Value = 30
$ws.Range("I9").Value = -10
$ws.Range("J9").Value = "LOW THREAT"
$ws.Range("K9").Value = "SAR"

# Row 10
$ws.Range("A10").Value = "06-FEB-26"
$ws.Range("B10").Value = "SM-438"
$ws.Range("C10").Value = "flynas XY-894"
$ws.Range("D10").Value = 709
$ws.Range("E10").Value = 830
$ws.Range("F10").Value = -121
$ws.Range("G10").Value = 40
$ws.Range("H10").Value = 30
$ws.Range("I10").Value = -10
$ws.Range("J10").Value = "LOW THREAT"
$ws.Range("K10").Value = "SAR"

# Row 11
$ws.Range("A11").Value = "06-FEB-26"
$ws.Range("B11").Value = "SM-438"
$ws.Range("C11").Value = "EgyptAir MS-684"
$ws.Range("D11").Value = 856
$ws.Range("E11").Value = 830
$ws.Range("F11").Value = 26
$ws.Range("G11").Value = 46
$ws.Range("H11").Value = 30
$ws.Range("I11").Value = -16
$ws.Range("J11").Value = "LOW THREAT"
$ws.Range("K11").Value = "SAR"

# Row 12
$ws.Range("A12").Value = "09-FEB-26"
$ws.Range("B12").Value = "SM-438"
$ws.Range("C12").Value = "Nile Air NP-232"
$ws.Range("D12").Value = 563
$ws.Range("E12").Value = 830
$ws.Range("F12").Value = -267
$ws.Range("G12").Value = 30
$ws.Range("H12").Value = 30
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = "LOW THREAT"
$ws.Range("K12").Value = "SAR"

# Row 13
$ws.Range("A13").Value = "09-FEB-26"
$ws.Range("B13").Value = "SM-438"
$ws.Range("C13").Value = "flyadeal F3-911"
$ws.Range("D13").Value = 577
$ws.Range("E13").Value = 830
$ws.Range("F13").Value = -253
$ws.Range("G13").Value = 30
$ws.Range("H13").Value = 30
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = "LOW THREAT"
$ws.Range("K13").Value = "SAR"

# Row 14
$ws.Range("A14").Value = "09-FEB-26"
$ws.Range("B14").Value = "SM-438"
$ws.Range("C14").Value = "flynas XY-854"
$ws.Range("D14").Value = 599
$ws.Range("E14").Value = 830
$ws.Range("F14").Value = -231
$ws.Range("G14").Value = 40
$ws.Range("H14").Value = 30
$ws.Range("I14").Value = -10
$ws.Range("J14").Value = "LOW THREAT"
$ws.Range("K14").Value = "SAR"

# Row 15
$ws.Range("A15").Value = "09-FEB-26"
$ws.Range("B15").Value = "SM-438"
$ws.Range("C15").Value = "flynas XY-894"
$ws.Range("D15").Value = 629
$ws.Range("E15").Value = 830
$ws.Range("F15").Value = -201
$ws.Range("G15").Value = 40
$ws.Range("H15").Value = 30
$ws.Range("I15").Value = -10
$ws.Range("J15").Value = "LOW THREAT"
$ws.Range("K15").Value = "SAR"

# Row 16
$ws.Range("A16").Value = "13-FEB-26"
$ws.Range("B16").Value = "SM-438"
$ws.Range("C16").Value = "flynas XY-894"
$ws.Range("D16").Value = 669
$ws.Range("E16").Value = 830
$ws.Range("F16").Value = -161
$ws.Range("G16").Value = 40
$ws.Range("H16").Value = 30
$ws.Range("I16").Value = -10
$ws.Range("J16").Value = "LOW THREAT"
$ws.Range("K16").Value = "SAR"

# Row 17
$ws.Range("A17").Value = "13-FEB-26"
$ws.Range("B17").Value = "SM-438"
$ws.Range("C17").Value = "EgyptAir MS-684"
$ws.Range("D17").Value = 848
$ws.Range("E17").Value = 830
$ws.Range("F17").Value = 18
$ws.Range("G17").Value = 46
$ws.Range("H17").Value = 30
$ws.Range("I17").Value = -16
$ws.Range("J17").Value = "LOW THREAT"
$ws.Range("K17").Value = "SAR"

# Row 18
$ws.Range("A18").Value = "16-FEB-26"
$ws.Range("B18").Value = "SM-438"
$ws.Range("C18").Value = "flyadeal F3-911"
$ws.Range("D18").Value = 547
$ws.Range("E18").Value = 826
$ws.Range("F18").Value = -279
$ws.Range("G18").Value = 30
$ws.Range("H18").Value = 30
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = "LOW THREAT"
$ws.Range("K18").Value = "SAR"

# Row 19
$ws.Range("A19").Value = "16-FEB-26"
$ws.Range("B19").Value = "SM-438"
$ws.Range("C19").Value = "Nile Air NP-232"
$ws.Range("D19").Value = 563
$ws.Range("E19").Value = 826
$ws.Range("F19").Value = -263
$ws.Range("G19").Value = 30
$ws.Range("H19").Value = 30
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = "LOW THREAT"
$ws.Range("K19").Value = "SAR"

# Row 20
$ws.Range("A20").Value = "16-FEB-26"
$ws.Range("B20").Value = "SM-438"
$ws.Range("C20").Value = "flynas XY-894"
$ws.Range("D20").Value = 599
$ws.Range("E20").Value = 826
$ws.Range("F20").Value = -227
$ws.Range("G20").Value = 40
$ws.Range("H20").Value = 30
$ws.Range("I20").Value = -10
$ws.Range("J20").Value = "LOW THREAT"
$ws.Range("K20").Value = "SAR"

# Row 21
$ws.Range("A21").Value = "16-FEB-26"
$ws.Range("B21").Value = "SM-438"
$ws.Range("C21").Value = "flynas XY-854"
$ws.Range("D21").Value = 599
$ws.Range("E21").Value = 826
$ws.Range("F21").Value = -227
$ws.Range("G21").Value = 40
$ws.Range("H21").Value = 30
$ws.Range("I21").Value = -10
$ws.Range("J21").Value = "LOW THREAT"
$ws.Range("K21").Value = "SAR"

# Row 22
$ws.Range("A22").Value = "16-FEB-26"
$ws.Range("B22").Value = "SM-438"
$ws.Range("C22").Value = "EgyptAir MS-682"
$ws.Range("D22").Value = 776
$ws.Range("E22").Value = 826
$ws.Range("F22").Value = -50
$ws.Range("G22").Value = 46
$ws.Range("H22").Value = 30
$ws.Range("I22").Value = -16
$ws.Range("J22").Value = "LOW THREAT"
$ws.Range("K22").Value = "SAR"

# Row 23
$ws.Range("A23").Value = "20-FEB-26"
$ws.Range("B23").Value = "SM-438"
$ws.Range("C23").Value = "flynas XY-894"
$ws.Range("D23").Value = 569
$ws.Range("E23").Value = 686
$ws.Range("F23").Value = -117
$ws.Range("G23").Value = 40
$ws.Range("H23").Value = 30
$ws.Range("I23").Value = -10
$ws.Range("J23").Value = "LOW THREAT"
$ws.Range("K23").Value = "SAR"

# Row 24
$ws.Range("A24").Value = "23-FEB-26"
$ws.Range("B24").Value = "SM-438"
$ws.Range("C24").Value = "Nile Air NP-232"
$ws.Range("D24").Value = 418
$ws.Range("E24").Value = 686
$ws.Range("F24").Value = -268
$ws.Range("G24").Value = 30
$ws.Range("H24").Value = 30
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = "LOW THREAT"
$ws.Range("K24").Value = "SAR"

# Row 25
$ws.Range("A25").Value = "23-FEB-26"
$ws.Range("B25").Value = "SM-438"
$ws.Range("C25").Value = "flynas XY-894"
$ws.Range("D25").Value = 529
$ws.Range("E25").Value = 686
$ws.Range("F25").Value = -157
$ws.Range("G25").Value = 40
$ws.Range("H25").Value = 30
$ws.Range("I25").Value = -10
$ws.Range("J25").Value = "LOW THREAT"
$ws.Range("K25").Value = "SAR"

# Row 26
$ws.Range("A26").Value = "23-FEB-26"
$ws.Range("B26").Value = "SM-438"
$ws.Range("C26").Value = "flynas XY-854"
$ws.Range("D26").Value = 529
$ws.Range("E26").Value = 686
$ws.Range("F26").Value = -157
$ws.Range("G26").Value = 40
$ws.Range("H26").Value = 30
$ws.Range("I26").Value = -10
$ws.Range("J26").Value = "LOW THREAT"
$ws.Range("K26").Value = "SAR"
